$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.167195999999999
$ws.Range("H2").Value = 24.501588
$ws.Range("I2").Value = 0.0776491163898855
$ws.Range("J2").Value = 0.07764911638988552
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 15.67208412997466
$ws.Range("R2").Value = 141.048757169772
$ws.Range("S2").Value = 0.000506584323279491
$ws.Range("T2").Value = 0.0005065843232794911
$ws.Range("G3").Value = 8.167195999999999
$ws.Range("H3").Value = 24.501588
$ws.Range("I3").Value = 0.0776491163898855
$ws.Range("J3").Value = 0.07764911638988552
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 1480.617824544035
$ws.Range("R3").Value = 13325.56042089631
$ws.Range("S3").Value = 0.04785947883266018
$ws.Range("T3").Value = 0.04785947883266019
$ws.Range("G4").Value = 8.167195999999999
$ws.Range("H4").Value = 24.501588
$ws.Range("I4").Value = 0.0776491163898855
$ws.Range("J4").Value = 0.07764911638988552
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 237.2380332901106
$ws.Range("R4").Value = 2135.142299610996
$ws.Range("S4").Value = 0.007668480308918706
$ws.Range("T4").Value = 0.007668480308918708
$ws.Range("G5").Value = 8.167195999999999
$ws.Range("H5").Value = 24.501588
$ws.Range("I5").Value = 0.0776491163898855
$ws.Range("J5").Value = 0.07764911638988552
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 668.6851324603786
$ws.Range("R5").Value = 6018.166192143407
$ws.Range("S5").Value = 0.02161457292502713
$ws.Range("T5").Value = 0.02161457292502713
$ws.Range("I6").Value = 0.6024007145055783
$ws.Range("J6").Value = 0.6024007145055783
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 121.5838005198219
$ws.Range("R6").Value = 1094.254204678397
$ws.Range("S6").Value = 0.0039300738049434
$ws.Range("T6").Value = 0.0039300738049434
$ws.Range("I7").Value = 0.6024007145055783
$ws.Range("J7").Value = 0.6024007145055783
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("S7").Value = 0.3712931400261824
$ws.Range("T7").Value = 0.3712931400261824
$ws.Range("I8").Value = 0.6024007145055783
$ws.Range("J8").Value = 0.6024007145055783
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 1840.489208457708
$ws.Range("R8").Value = 16564.40287611937
$ws.Range("S8").Value = 0.05949195859576217
$ws.Range("T8").Value = 0.05949195859576218
$ws.Range("I9").Value = 0.6024007145055783
$ws.Range("J9").Value = 0.6024007145055783
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 5187.649522639767
$ws.Range("R9").Value = 46688.84570375791
$ws.Range("S9").Value = 0.1676855420786903
$ws.Range("T9").Value = 0.1676855420786903
$ws.Range("G10").Value = 25.97600266666667
$ws.Range("H10").Value = 77.92800800000001
$ws.Range("I10").Value = 0.2469652564243563
$ws.Range("J10").Value = 0.2469652564243563
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 49.84551603175022
$ws.Range("R10").Value = 448.609644285752
$ws.Range("S10").Value = 0.001611206065386406
$ws.Range("T10").Value = 0.001611206065386406
$ws.Range("G11").Value = 25.97600266666667
$ws.Range("H11").Value = 77.92800800000001
$ws.Range("I11").Value = 0.2469652564243563
$ws.Range("J11").Value = 0.2469652564243563
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 4709.147736710378
$ws.Range("R11").Value = 42382.3296303934
$ws.Range("S11").Value = 0.1522184541404979
$ws.Range("T11").Value = 0.1522184541404979
$ws.Range("G12").Value = 25.97600266666667
$ws.Range("H12").Value = 77.92800800000001
$ws.Range("I12").Value = 0.2469652564243563
$ws.Range("J12").Value = 0.2469652564243563
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 754.5424139911263
$ws.Range("R12").Value = 6790.881725920137
$ws.Range("S12").Value = 0.0243898230131557
$ws.Range("T12").Value = 0.02438982301315571
$ws.Range("G13").Value = 25.97600266666667
$ws.Range("H13").Value = 77.92800800000001
$ws.Range("I13").Value = 0.2469652564243563
$ws.Range("J13").Value = 0.2469652564243563
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 2126.772368870681
$ws.Range("R13").Value = 19140.95131983613
$ws.Range("S13").Value = 0.06874577320531623
$ws.Range("T13").Value = 0.06874577320531625
$ws.Range("G14").Value = 7.676611333333334
$ws.Range("H14").Value = 23.029834
$ws.Range("I14").Value = 0.07298491268017987
$ws.Range("J14").Value = 0.07298491268017987
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 14.73069810607178
$ws.Range("R14").Value = 132.576282954646
$ws.Range("S14").Value = 0.0004761549688995267
$ws.Range("T14").Value = 0.0004761549688995267
$ws.Range("G15").Value = 7.676611333333334
$ws.Range("H15").Value = 23.029834
$ws.Range("I15").Value = 0.07298491268017987
$ws.Range("J15").Value = 0.07298491268017987
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 1391.680519511235
$ws.Range("R15").Value = 12525.12467560112
$ws.Range("S15").Value = 0.04498467008924801
$ws.Range("T15").Value = 0.04498467008924801
$ws.Range("G16").Value = 7.676611333333334
$ws.Range("H16").Value = 23.029834
$ws.Range("I16").Value = 0.07298491268017987
$ws.Range("J16").Value = 0.07298491268017987
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 222.9876906410198
$ws.Range("R16").Value = 2006.889215769178
$ws.Range("S16").Value = 0.007207852346005759
$ws.Range("T16").Value = 0.00720785234600576
$ws.Range("G17").Value = 7.676611333333334
$ws.Range("H17").Value = 23.029834
$ws.Range("I17").Value = 0.07298491268017987
$ws.Range("J17").Value = 0.07298491268017987
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 628.5187555529271
$ws.Range("R17").Value = 5656.668799976344
$ws.Range("S17").Value = 0.02031623527602657
$ws.Range("T17").Value = 0.02031623527602657
